# Daily attendance processing - 2025-11-14 12:39:54
# Normalize the "Recorded By" (column G) lists so that the "System" entry
# is listed last instead of first, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $parts = $val.Split(",")
    if ($parts.Count -gt 1) {
        $first = $parts[0].Trim()
        if ($first.ToLower() -eq "system") {
            $rest = @()
            for ($i = 1; $i -lt $parts.Count; $i++) {
                $rest += $parts[$i].Trim()
            }
            $rest += $first
            $newVal = [string]::Join(", ", $rest)
            $cell.Value = $newVal
            $changed++
        }
    }
}

Write-Host ("Reordered 'System' entry in " + $changed + " 'Recorded By' cells.")
